$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 142
$ws.Range("B142").Value = 63902
$ws.Range("E142").Value = 34.04
$ws.Range("F142").Value = 2
$ws.Range("G142").Value = 64.04000000000001

# Row 143
$ws.Range("B143").Value = 48654
$ws.Range("E143").Value = 38.26
$ws.Range("F143").Value = -1
$ws.Range("G143").Value = -32.02

# Row 176
$ws.Range("B176").Value = 57552
$ws.Range("E176").Value = 136.86
$ws.Range("F176").Value = -5
$ws.Range("G176").Value = -603.45

# Row 177
$ws.Range("B177").Value = 64329
$ws.Range("E177").Value = 128.32
$ws.Range("F177").Value = 6
$ws.Range("G177").Value = 724.14

# Row 256
$ws.Range("B256").Value = 64979
$ws.Range("E256").Value = 314.41
$ws.Range("F256").Value = 82
$ws.Range("G256").Value = 24251.5

# Row 257
$ws.Range("B257").Value = 48719
$ws.Range("E257").Value = 353.35
$ws.Range("F257").Value = -81
$ws.Range("G257").Value = -23955.75

# Row 271
$ws.Range("B271").Value = 64973
$ws.Range("E271").Value = 35.4
$ws.Range("F271").Value = 150
$ws.Range("G271").Value = 4995

# Row 272
$ws.Range("B272").Value = 48706
$ws.Range("E272").Value = 39.8
$ws.Range("F272").Value = -144
$ws.Range("G272").Value = -4795.2

# Row 308
$ws.Range("B308").Value = 63565
$ws.Range("D308").Value = 102.71
$ws.Range("E308").Value = 109.19
$ws.Range("F308").Value = 60
$ws.Range("G308").Value = 6162.6

# Row 309
$ws.Range("B309").Value = 57077
$ws.Range("D309").Value = 93.08
$ws.Range("E309").Value = 111.2
$ws.Range("F309").Value = 1
$ws.Range("G309").Value = 93.08

# Row 342
$ws.Range("B342").Value = 63571
$ws.Range("E342").Value = 152.53
$ws.Range("F342").Value = 29
$ws.Range("G342").Value = 4160.92

# Row 343
$ws.Range("B343").Value = 57802
$ws.Range("E343").Value = 162.71
$ws.Range("F343").Value = -79
$ws.Range("G343").Value = -11334.92

# Row 347
$ws.Range("B347").Value = 63510
$ws.Range("E347").Value = 50.66
$ws.Range("F347").Value = 167
$ws.Range("G347").Value = 7955.88

# Row 348
$ws.Range("B348").Value = 55356
$ws.Range("E348").Value = 54.04
$ws.Range("F348").Value = -158
$ws.Range("G348").Value = -7527.12

# Row 367
$ws.Range("B367").Value = 61605
$ws.Range("E367").Value = 133.78
$ws.Range("F367").Value = -13
$ws.Range("G367").Value = -1455.48

# Row 368
$ws.Range("B368").Value = 63563
$ws.Range("E368").Value = 119.04
$ws.Range("F368").Value = 15
$ws.Range("G368").Value = 1679.4

# Row 411
$ws.Range("B411").Value = 57856
$ws.Range("F411").Value = 2
$ws.Range("G411").Value = 342.66

# Row 412
$ws.Range("B412").Value = 63007
$ws.Range("F412").Value = 984
$ws.Range("G412").Value = 168588.72

# Row 423
$ws.Range("B423").Value = 63102
$ws.Range("C423").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F423").Value = 36
$ws.Range("G423").Value = 2140.92

# Row 424
$ws.Range("B424").Value = 53082
$ws.Range("C424").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F424").Value = 1
$ws.Range("G424").Value = 59.47

# Row 449
$ws.Range("B449").Value = 63681
$ws.Range("E449").Value = 23.84
$ws.Range("F449").Value = 65
$ws.Range("G449").Value = 1457.3

# Row 450
$ws.Range("B450").Value = 31930
$ws.Range("E450").Value = 26.8
$ws.Range("F450").Value = -62
$ws.Range("G450").Value = -1390.04

# Row 528
$ws.Range("B528").Value = 58047
$ws.Range("D528").Value = 105.54
$ws.Range("E528").Value = 126.1
$ws.Range("F528").Value = 54
$ws.Range("G528").Value = 5699.16

# Row 529
$ws.Range("B529").Value = 47097
$ws.Range("D529").Value = 112.28
$ws.Range("E529").Value = 134.16
$ws.Range("F529").Value = 15
$ws.Range("G529").Value = 1684.2

# Row 571
$ws.Range("B571").Value = 53757
$ws.Range("E571").Value = 16.08
$ws.Range("F571").Value = -159
$ws.Range("G571").Value = -2138.55

# Row 572
$ws.Range("B572").Value = 65069
$ws.Range("E572").Value = 14.3
$ws.Range("F572").Value = 172
$ws.Range("G572").Value = 2313.4

# Row 578
$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28

# Row 579
$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2

# Row 582
$ws.Range("B582").Value = 45706
$ws.Range("E582").Value = 23.58
$ws.Range("F582").Value = -202
$ws.Range("G582").Value = -3985.46

# Row 583
$ws.Range("B583").Value = 64922
$ws.Range("E583").Value = 20.98
$ws.Range("F583").Value = 207
$ws.Range("G583").Value = 4084.11

# Row 585
$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9

# Row 586
$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68

# Row 591
$ws.Range("B591").Value = 64925
$ws.Range("E591").Value = 13.97
$ws.Range("F591").Value = 302
$ws.Range("G591").Value = 3971.3

# Row 592
$ws.Range("B592").Value = 45709
$ws.Range("E592").Value = 15.69
$ws.Range("F592").Value = -300
$ws.Range("G592").Value = -3945
